$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.478.60'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.87%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.983.89'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.47'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.37%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.78%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.592'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.93%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.74'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.28%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.71%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.37%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.453.17'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.34%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.45'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.35%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.79'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.06%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.990.53'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.57%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '11.10'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.94%  '

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.38%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '51.462.54'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.92%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.09'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.16%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.61'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.71%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.76%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.52'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.73%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.44'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.97%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.55%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.85'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.64%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.36'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.05%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.02%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.168'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.62%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.09'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.03%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.110'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.29%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.32'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.25%  '

$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '51.59'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.65%  '

$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '34.61'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.95%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.06'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.29%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0440'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.34%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.03%  '

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.59%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.78'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.90%  '

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.73%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.08%  '

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.72%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '124.84'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.63%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.74'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +12.27%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.69%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.01%  '

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.88%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.271'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.032.77'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.21%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.540'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +16.88%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0331'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.47%  '
